{"js": "// Add a new \"Phase 1 Senior Design Report\" bullet item to the README's\n// bulleted file list, right after the existing \"Bill Of Materials Appendix\"\n// item (the last paragraph in the document body).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The new list item goes after the very last paragraph in the document\n// (the \"Bill Of Materials Appendix\" bullet).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Build the new bullet as a Flat-OPC OOXML fragment so the bold \"label\" run\n// and the plain-text \"description\" run come out exactly as two distinct\n// <w:r> runs (one bold, one with no run properties at all) instead of Word\n// leaving behind an empty <w:rPr/> on the second run.\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"ListParagraph\"/>' +\n  '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n  '<w:rPr><w:b/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t>Phase 1 Senior Design Report</w:t></w:r>' +\n  '<w:r><w:t>: Word document containing the report written for phase 1 of this project.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// A caret positioned immediately after the last paragraph (i.e. past its\n// paragraph mark), so the inserted OOXML becomes a brand-new sibling\n// paragraph rather than splicing into the existing one.\nconst insertionPoint = lastParagraph.getRange(Word.RangeLocation.after);\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new \"Phase 1 Senior Design Report\" bullet item to the README's\n# bulleted file list, right after the existing \"Bill Of Materials Appendix\"\n# item (the last paragraph in the document).\n\n$d = $word.ActiveDocument\n\n# Caret at the very end of the document's story (after the last paragraph's\n# mark, i.e. past the \"Bill Of Materials Appendix\" bullet and its bookmark).\n$endRange = $d.Content\n$endRange.Collapse(0)  # wdCollapseEnd\n\n# Flat-OPC OOXML fragment for the new bullet: a bold \"label\" run followed by\n# a plain run, so the plain run has no <w:rPr> at all (matching how Word\n# itself emits an un-formatted run) instead of an explicit/empty one.\n$xml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:b/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n              </w:rPr>\n              <w:t>Phase 1 Senior Design Report</w:t>\n            </w:r>\n            <w:r>\n              <w:t>: Word document containing the report written for phase 1 of this project.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$endRange.InsertXML($xml)\n"}
